# Update the 15 lattice-multiplication exercise cells to the new set of
# problems. The table layout (5 rows x 3 columns) is unchanged; only the
# text content of each cell (the multiplication expression, the two
# multiplier digits, the divider, and the two placeholder rows) changes.
#
# Each cell holds a single run containing five <w:t> text nodes separated
# by manual line breaks (<w:br/>). We rebuild that run's paragraph via
# Range.InsertXML so we can control w:rPr (sz 32) and xml:space exactly,
# instead of relying on Find/Replace (which would upset whitespace
# handling across the whole document).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New content, in row-major reading order (row 1 col 1, row 1 col 2, ...),
# matching the 5x3 table exactly.
$cellData = @(
    @('79 x 15', '  1    5', '  ----', '7|    |', '9|    |'),
    @('34 x 57', '  5    7', '  ----', '3|    |', '4|    |'),
    @('94 x 59', '  5    9', '  ----', '9|    |', '4|    |'),
    @('35 x 52', '  5    2', '  ----', '3|    |', '5|    |'),
    @('95 x 93', '  9    3', '  ----', '9|    |', '5|    |'),
    @('98 x 69', '  6    9', '  ----', '9|    |', '8|    |'),
    @('86 x 53', '  5    3', '  ----', '8|    |', '6|    |'),
    @('53 x 97', '  9    7', '  ----', '5|    |', '3|    |'),
    @('73 x 58', '  5    8', '  ----', '7|    |', '3|    |'),
    @('44 x 22', '  2    2', '  ----', '4|    |', '4|    |'),
    @('46 x 51', '  5    1', '  ----', '4|    |', '6|    |'),
    @('48 x 36', '  3    6', '  ----', '4|    |', '8|    |'),
    @('34 x 14', '  1    4', '  ----', '3|    |', '4|    |'),
    @('12 x 99', '  9    9', '  ----', '1|    |', '2|    |'),
    @('66 x 80', '  8    0', '  ----', '6|    |', '6|    |')
)

$cols = 3
for ($idx = 0; $idx -lt $cellData.Length; $idx++) {
    $row = [int]([math]::Floor($idx / $cols)) + 1
    $col = ($idx % $cols) + 1
    $lines = $cellData[$idx]

    $runXml = ""
    for ($i = 0; $i -lt $lines.Length; $i++) {
        $line = $lines[$i]
        if ($i -gt 0) {
            $runXml += "<w:br/>"
        }
        if ($line -ne $line.Trim()) {
            $runXml += '<w:t xml:space="preserve">' + $line + '</w:t>'
        } else {
            $runXml += '<w:t>' + $line + '</w:t>'
        }
    }

    $pXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + $runXml + '</w:r></w:p>'

    $cell = $t.Cell($row, $col)
    [void]$cell.Range.InsertXML($pXml)
}

Write-Host "Updated" $cellData.Length "cells"
